# Update countries & provincias Spain
# This script applies the data refresh captured in the commit:
#  - Update the "last refreshed" timestamp
#  - Re-rank Japon (overtakes Malasia) and Mexico (overtakes Tailandia)
#    by shifting the country labels in the affected row block
#  - Refresh "Nuevos casos" (C) / "Muertes hoy" (G) to 0 for this snapshot,
#    and push updated totals for the handful of countries with new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 03:22"

# 2) Re-rank countries: Japon moves above Malasia, Mexico moves above Tailandia.
#    Only the country-name labels in column A need to move; the row's own
#    numeric data is updated separately below.
$countryUpdates = @(
    @{ Row = 33; Name = "Japon" },
    @{ Row = 34; Name = "Malasia" },
    @{ Row = 35; Name = "Pakistan" },
    @{ Row = 36; Name = "Ecuador" },
    @{ Row = 37; Name = "Filipinas" },
    @{ Row = 42; Name = "Mexico" },
    @{ Row = 43; Name = "Tailandia" },
    @{ Row = 44; Name = "Serbia" },
    @{ Row = 45; Name = "Finlandia" }
)

foreach ($u in $countryUpdates) {
    $addr = "A" + $u.Row
    $ws.Range($addr).Value = $u.Name
}

# 3) Apply the numeric updates for every affected row/column
$rowUpdates = @(
    @{ Row = 4; C=0; G=0 },
    @{ Row = 5; C=0; G=0 },
    @{ Row = 6; C=0; G=0 },
    @{ Row = 7; B=103375; C=0; E=72865; G=0 },
    @{ Row = 8; C=0; G=0 },
    @{ Row = 9; B=81740; C=32; D=77167; E=1242; F=211 },
    @{ Row = 10; C=0; G=0 },
    @{ Row = 11; C=0; G=0 },
    @{ Row = 12; C=0; G=0 },
    @{ Row = 13; C=0; G=0 },
    @{ Row = 14; C=0; G=0 },
    @{ Row = 15; C=0; G=0 },
    @{ Row = 16; C=0; G=0 },
    @{ Row = 17; C=0; G=0 },
    @{ Row = 18; B=12232; C=49; E=11540; G=1; H=565 },
    @{ Row = 19; C=0; G=0 },
    @{ Row = 20; C=0; G=0 },
    @{ Row = 21; C=0; G=0 },
    @{ Row = 22; C=0; G=0 },
    @{ Row = 23; C=0; G=0 },
    @{ Row = 24; C=0; G=0 },
    @{ Row = 25; C=0; G=0 },
    @{ Row = 26; C=0; G=0 },
    @{ Row = 27; C=0; G=0 },
    @{ Row = 28; C=0; G=0 },
    @{ Row = 29; C=0; G=0 },
    @{ Row = 30; C=0; G=0 },
    @{ Row = 31; C=0; G=0 },
    @{ Row = 32; C=0; G=0 },
    @{ Row = 33; B=3906; C=0; D=592; E=3222; F=64; G=0; H=92 },
    @{ Row = 34; B=3793; C=0; D=1241; E=2490; F=102; G=0; H=62 },
    @{ Row = 35; B=3766; C=0; D=259; E=3454; F=17; G=0; H=53 },
    @{ Row = 36; B=3747; C=0; D=100; E=3456; F=156; G=0; H=191 },
    @{ Row = 37; B=3660; D=73; E=3424; F=1; H=163 },
    @{ Row = 38; C=0; G=0 },
    @{ Row = 39; C=0; G=0 },
    @{ Row = 40; C=0; G=0 },
    @{ Row = 41; C=0; G=0 },
    @{ Row = 42; B=2439; C=296; D=633; E=1681; F=419; G=31; H=125 },
    @{ Row = 43; B=2220; C=0; D=793; E=1401; F=23; G=0; H=26 },
    @{ Row = 44; B=2200; C=0; D=118; E=2024; F=101; H=58 },
    @{ Row = 45; B=2176; C=0; D=300; E=1849; F=81; G=0; H=27 },
    @{ Row = 46; C=0; G=0 },
    @{ Row = 47; C=0; G=0 },
    @{ Row = 48; C=0 },
    @{ Row = 49; C=0; G=0 },
    @{ Row = 50; C=0; G=0 },
    @{ Row = 51; C=0; G=0 },
    @{ Row = 52; C=0; G=0 },
    @{ Row = 53; C=0; G=0 },
    @{ Row = 54; C=0; G=0 },
    @{ Row = 55; C=0; G=0 },
    @{ Row = 56; C=0 },
    @{ Row = 57; C=0; G=0 },
    @{ Row = 58; C=0; G=0 },
    @{ Row = 59; C=0; G=0 },
    @{ Row = 60; C=0; G=0 },
    @{ Row = 61; C=0; G=0 },
    @{ Row = 62; C=0 },
    @{ Row = 63; C=0; G=0 },
    @{ Row = 64; C=0; G=0 },
    @{ Row = 65; C=0; G=0 },
    @{ Row = 66; C=0 },
    @{ Row = 67; C=0; G=0 },
    @{ Row = 68; C=0; G=0 },
    @{ Row = 69; C=0 },
    @{ Row = 70; C=0; G=0 },
    @{ Row = 72; C=0; G=0 },
    @{ Row = 73; C=0; G=0 },
    @{ Row = 74; C=0 },
    @{ Row = 75; C=0 },
    @{ Row = 76; C=0 },
    @{ Row = 77; C=0 },
    @{ Row = 78; C=0 },
    @{ Row = 79; C=0; G=0 },
    @{ Row = 80; C=0; G=0 },
    @{ Row = 81; C=0 },
    @{ Row = 82; C=0; G=0 },
    @{ Row = 83; C=0; G=0 },
    @{ Row = 84; C=0; G=0 },
    @{ Row = 85; C=0 },
    @{ Row = 86; C=0 },
    @{ Row = 87; C=0 },
    @{ Row = 88; C=0 },
    @{ Row = 89; C=0; G=0 },
    @{ Row = 90; C=0 },
    @{ Row = 91; C=0; G=0 },
    @{ Row = 92; C=0; G=0 },
    @{ Row = 93; C=0; G=0 },
    @{ Row = 94; C=0 },
    @{ Row = 95; C=0; G=0 },
    @{ Row = 96; C=0 },
    @{ Row = 97; C=0 },
    @{ Row = 98; C=0 },
    @{ Row = 99; C=0 },
    @{ Row = 100; C=0 },
    @{ Row = 101; C=0 },
    @{ Row = 102; C=0 },
    @{ Row = 103; C=0 },
    @{ Row = 104; C=0 },
    @{ Row = 105; C=0 },
    @{ Row = 106; C=0 },
    @{ Row = 107; C=0 },
    @{ Row = 108; C=0; G=0 },
    @{ Row = 110; C=0 },
    @{ Row = 111; C=0; G=0 },
    @{ Row = 112; C=0 },
    @{ Row = 113; C=0 },
    @{ Row = 114; C=0 },
    @{ Row = 115; C=0 },
    @{ Row = 116; C=0 },
    @{ Row = 117; C=0; G=0 },
    @{ Row = 118; C=0 },
    @{ Row = 119; C=0 },
    @{ Row = 120; C=0 },
    @{ Row = 122; C=0 },
    @{ Row = 123; C=0; G=0 },
    @{ Row = 125; C=0; G=0 },
    @{ Row = 126; C=0 },
    @{ Row = 127; C=0 },
    @{ Row = 128; C=0; G=0 },
    @{ Row = 129; C=0 },
    @{ Row = 130; C=0 },
    @{ Row = 131; C=0 },
    @{ Row = 133; C=0 },
    @{ Row = 134; C=0 },
    @{ Row = 136; C=0; G=0 },
    @{ Row = 137; C=0; G=0 },
    @{ Row = 138; C=0 },
    @{ Row = 139; C=0 },
    @{ Row = 141; C=0 },
    @{ Row = 143; C=0 },
    @{ Row = 145; C=0 },
    @{ Row = 149; C=0; G=0 },
    @{ Row = 150; C=0; G=0 },
    @{ Row = 151; C=0 },
    @{ Row = 154; C=0 },
    @{ Row = 155; C=0 },
    @{ Row = 156; C=0; G=0 },
    @{ Row = 157; C=0 },
    @{ Row = 158; C=0 },
    @{ Row = 159; C=0 },
    @{ Row = 160; C=0 },
    @{ Row = 161; C=0 },
    @{ Row = 169; C=0 },
    @{ Row = 171; C=0 },
    @{ Row = 172; C=0 },
    @{ Row = 173; C=0 },
    @{ Row = 175; C=0 },
    @{ Row = 176; C=0 },
    @{ Row = 178; C=0 },
    @{ Row = 180; C=0 },
    @{ Row = 183; C=0 },
    @{ Row = 186; C=0 },
    @{ Row = 191; C=0 },
    @{ Row = 194; C=0; G=0 },
    @{ Row = 202; C=0 },
    @{ Row = 204; C=0 },
    @{ Row = 212; C=0 }
)

$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($u in $rowUpdates) {
    $r = $u.Row
    foreach ($col in $cols) {
        if ($u.ContainsKey($col)) {
            $addr = "$col$r"
            $ws.Range($addr).Value = $u[$col]
        }
    }
}
